$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D as Text (it holds decimal-looking strings, e.g. "0.9998", "1.847.07")
# so COM does not silently coerce them into numeric Variants on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.414.21'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.849.94'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '240.83'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '0.6300'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('D8').Value = '0.07684'
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('D9').Value = '0.2942'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '24.53'
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('D11').Value = '0.07749'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').Value = '1.847.07'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '0.00001085'
$ws.Range('E14').Value = '  +8.14%  '
$ws.Range('D15').Value = '0.6814'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '83.79'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '2.097.48'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '6.147'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '29.421.40'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '229.39'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '7.451'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '157.09'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('D27').Value = '8.393'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '17.68'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').Value = '1.319'
$ws.Range('E29').Value = '  +4.20%  '
$ws.Range('D30').Value = '1.466'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').Value = '0.05697'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '4.114'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = '4.048'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').Value = '1.852'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '0.7077'
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '2.779'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = '0.01799'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '1.224.68'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').Value = '6.441'
$ws.Range('E41').Value = '  +4.37%  '
$ws.Range('D42').Value = '0.9098'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '2.006.58'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '101.55'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '66.15'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000121'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '7.136'
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '0.4024'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.041'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.690'
$ws.Range('E51').Value = '  +0.33%  '
